# Scheduled-runner refresh of market-price-derived leve profit figures
# (columns H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 215.5
$ws.Range("J9").Value = 200
$ws.Range("L9").Value = 200
$ws.Range("N9").Value = -538
$ws.Range("H11").Value = 74.42856999999999
$ws.Range("I11").Value = 74.42856999999999
$ws.Range("K11").Value = 74.42856999999999
$ws.Range("M11").Value = 65.57143000000001
$ws.Range("H15").Value = 1903.4773
$ws.Range("I15").Value = 1903.4773
$ws.Range("K15").Value = 5710.4319
$ws.Range("M15").Value = -5541.4319
$ws.Range("H17").Value = 2511.2856
$ws.Range("I17").Value = 996
$ws.Range("J17").Value = 2763.8333
$ws.Range("K17").Value = 2988
$ws.Range("L17").Value = 8291.499899999999
$ws.Range("M17").Value = -2820
$ws.Range("N17").Value = -8627.499899999999
$ws.Range("H28").Value = 7106.7144
$ws.Range("I28").Value = 1550
$ws.Range("K28").Value = 1550
$ws.Range("M28").Value = -1065
$ws.Range("H40").Value = 3568.6216
$ws.Range("I40").Value = 2720
$ws.Range("K40").Value = 2720
$ws.Range("M40").Value = -2545
$ws.Range("H51").Value = 2922.5
$ws.Range("J51").Value = 2922.5
$ws.Range("L51").Value = 2922.5
$ws.Range("N51").Value = -3890.5
$ws.Range("H92").Value = 225.5
$ws.Range("I92").Value = 210.6
$ws.Range("K92").Value = 210.6
$ws.Range("M92").Value = 1037.4
$ws.Range("H101").Value = 324.83334
$ws.Range("J101").Value = 350
$ws.Range("L101").Value = 1050
$ws.Range("N101").Value = -4294
$ws.Range("H107").Value = 406
$ws.Range("I107").Value = 118.333336
$ws.Range("J107").Value = 2995
$ws.Range("K107").Value = 118.333336
$ws.Range("L107").Value = 2995
$ws.Range("M107").Value = 1801.666664
$ws.Range("N107").Value = -6835
$ws.Range("H113").Value = 4999
$ws.Range("H127").Value = 785.75
$ws.Range("I127").Value = 764.3333
$ws.Range("J127").Value = 850
$ws.Range("K127").Value = 2292.9999
$ws.Range("L127").Value = 2550
$ws.Range("M127").Value = 2667.0001
$ws.Range("N127").Value = -12470
$ws.Range("H131").Value = 3783.125
$ws.Range("I131").Value = 3877.5
$ws.Range("K131").Value = 11632.5
$ws.Range("M131").Value = -6592.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8691
$ws.Range("J2").Value = 8756.5
$ws.Range("L2").Value = 8756.5
$ws.Range("N2").Value = -8982.5
$ws.Range("H4").Value = 385.2857
$ws.Range("I4").Value = 282.83334
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 282.83334
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -166.83334
$ws.Range("N4").Value = -1232
$ws.Range("H110").Value = 2291.0667
$ws.Range("I110").Value = 1755.125
$ws.Range("K110").Value = 1755.125
$ws.Range("M110").Value = 289.875
$ws.Range("H116").Value = 8691
$ws.Range("J116").Value = 8756.5
$ws.Range("L116").Value = 8756.5
$ws.Range("N116").Value = -13344.5
$ws.Range("H122").Value = 2957.1667
$ws.Range("I122").Value = 2798.6
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 8395.799999999999
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -5945.799999999999
$ws.Range("N122").Value = -16150
$ws.Range("H132").Value = 2189.8708
$ws.Range("I132").Value = 2277.2964
$ws.Range("K132").Value = 6831.889200000001
$ws.Range("M132").Value = -4301.889200000001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8691
$ws.Range("J3").Value = 8756.5
$ws.Range("L3").Value = 8756.5
$ws.Range("N3").Value = -8984.5
$ws.Range("H75").Value = 55000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 55000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 55000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -56872
$ws.Range("H78").Value = 55000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 55000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 165000
$ws.Range("M78").Value = -174360
$ws.Range("H99").Value = 1736.6428
$ws.Range("I99").Value = 1837.7273
$ws.Range("J99").Value = 1366
$ws.Range("K99").Value = 1837.7273
$ws.Range("L99").Value = 1366
$ws.Range("M99").Value = -339.7273
$ws.Range("N99").Value = -4362
$ws.Range("H105").Value = 2500
$ws.Range("I105").Value = 2500
$ws.Range("K105").Value = 2500
$ws.Range("M105").Value = -753
$ws.Range("H107").Value = 4964.8125
$ws.Range("I107").Value = 1048.5555
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 1048.5555
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 871.4445000000001
$ws.Range("N107").Value = -13840
$ws.Range("H134").Value = 2201.2307
$ws.Range("I134").Value = 1909.75
$ws.Range("K134").Value = 5729.25
$ws.Range("M134").Value = -3194.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1878.5714
$ws.Range("I16").Value = 1191.6666
$ws.Range("K16").Value = 1191.6666
$ws.Range("M16").Value = -904.6666
$ws.Range("H31").Value = 4675.731
$ws.Range("I31").Value = 3532.5557
$ws.Range("J31").Value = 5280.9414
$ws.Range("K31").Value = 3532.5557
$ws.Range("L31").Value = 5280.9414
$ws.Range("M31").Value = -3237.5557
$ws.Range("N31").Value = -5870.9414
$ws.Range("H34").Value = 4675.731
$ws.Range("I34").Value = 3532.5557
$ws.Range("J34").Value = 5280.9414
$ws.Range("K34").Value = 3532.5557
$ws.Range("L34").Value = 5280.9414
$ws.Range("M34").Value = -3330.5557
$ws.Range("N34").Value = -5684.9414
$ws.Range("H107").Value = 1127
$ws.Range("I107").Value = 914.4545000000001
$ws.Range("K107").Value = 914.4545000000001
$ws.Range("M107").Value = 1005.5455
$ws.Range("H113").Value = 1878.5714
$ws.Range("I113").Value = 1191.6666
$ws.Range("K113").Value = 1191.6666
$ws.Range("M113").Value = 978.3334
$ws.Range("H122").Value = 1129.3334
$ws.Range("I122").Value = 1055.2
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3165.6
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -715.6000000000004
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 4097.7144
$ws.Range("I132").Value = 3736.8
$ws.Range("K132").Value = 11210.4
$ws.Range("M132").Value = -8680.400000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1790.25
$ws.Range("I34").Value = 357
$ws.Range("J34").Value = 2268
$ws.Range("K34").Value = 1071
$ws.Range("L34").Value = 6804
$ws.Range("M34").Value = -987
$ws.Range("N34").Value = -6972
$ws.Range("H56").Value = 6665
$ws.Range("I56").Value = 6665
$ws.Range("K56").Value = 6665
$ws.Range("M56").Value = -6135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 25.444445
$ws.Range("I2").Value = 11.666667
$ws.Range("K2").Value = 11.666667
$ws.Range("M2").Value = 101.333333
$ws.Range("H80").Value = 2598.6667
$ws.Range("I80").Value = 2024.25
$ws.Range("K80").Value = 2024.25
$ws.Range("M80").Value = -1026.25
$ws.Range("H83").Value = 2598.6667
$ws.Range("I83").Value = 2024.25
$ws.Range("K83").Value = 10121.25
$ws.Range("M83").Value = -5129.25
$ws.Range("H113").Value = 7385.231
$ws.Range("I113").Value = 3599.6
$ws.Range("J113").Value = 9751.25
$ws.Range("K113").Value = 3599.6
$ws.Range("L113").Value = 9751.25
$ws.Range("M113").Value = -1429.6
$ws.Range("N113").Value = -14091.25
$ws.Range("H122").Value = 2426.75
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 63239.59
$ws.Range("I132").Value = 103002.5
$ws.Range("K132").Value = 309007.5
$ws.Range("M132").Value = -306477.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -939
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4696
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 906.2222
$ws.Range("I107").Value = 516.4
$ws.Range("J107").Value = 1393.5
$ws.Range("K107").Value = 1549.2
$ws.Range("L107").Value = 4180.5
$ws.Range("M107").Value = 370.8000000000002
$ws.Range("N107").Value = -8020.5
$ws.Range("H122").Value = 4488
$ws.Range("I122").Value = 1642.6666
$ws.Range("J122").Value = 7333.3335
$ws.Range("K122").Value = 4927.9998
$ws.Range("L122").Value = 22000.0005
$ws.Range("M122").Value = -2477.9998
$ws.Range("N122").Value = -26900.0005
$ws.Range("H132").Value = 2609.4443
$ws.Range("I132").Value = 2570.7144
$ws.Range("J132").Value = 2745
$ws.Range("K132").Value = 7712.1432
$ws.Range("L132").Value = 8235
$ws.Range("M132").Value = -5182.1432
$ws.Range("N132").Value = -13295
